# "using gains for all" -- insert two new metric columns (M_TotalTax, M_CorpTax)
# right after M_POP, pushing the existing GFA/IMF/OECD gain columns two slots
# to the right, and fix up one data value that changed alongside it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank columns at F:G -- everything from the old F column
# (GFA - Sales) onward shifts right to H:O.
$ws.Range("F1:G1").EntireColumn.Insert()

# New header labels for the inserted columns.
$ws.Range("F1").Value = "M_TotalTax"
$ws.Range("G1").Value = "M_CorpTax"

# New M_TotalTax / M_CorpTax data values for each region row.
$ws.Range("F2").Value = 5453258922371.839
$ws.Range("G2").Value = 901477052802.8717

$ws.Range("F3").Value = 7186470855993.515
$ws.Range("G3").Value = 480200693963.0637

$ws.Range("F4").Value = 1002696890625.219
$ws.Range("G4").Value = 134030001792.1409

$ws.Range("F5").Value = 241723618867.4443
$ws.Range("G5").Value = 38152875196.83905

$ws.Range("F6").Value = 5526698445364.336
$ws.Range("G6").Value = 408160442622.5067

$ws.Range("F7").Value = 457036063703.7742
$ws.Range("G7").Value = 4932505470.985653

$ws.Range("F8").Value = 114563677189.0651
$ws.Range("G8").Value = 24178295225.23782

# M_POP value for Sub - Saharan Africa (row 8, column E) was also recomputed.
$ws.Range("E8").Value = 366265684.25
